$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H29").Value = 2043.75
$ws_ALC.Range("J29").Value = 2625
$ws_ALC.Range("L29").Value = 7875
$ws_ALC.Range("N29").Value = -8437
$ws_ALC.Range("H33").Value = 261.2857
$ws_ALC.Range("I33").Value = 230.6875
$ws_ALC.Range("J33").Value = 359.2
$ws_ALC.Range("K33").Value = 230.6875
$ws_ALC.Range("L33").Value = 359.2
$ws_ALC.Range("M33").Value = -1.6875
$ws_ALC.Range("N33").Value = -817.2
$ws_ALC.Range("H38").Value = 1192.1818
$ws_ALC.Range("I38").Value = 232.8
$ws_ALC.Range("J38").Value = 1991.6666
$ws_ALC.Range("K38").Value = 698.4000000000001
$ws_ALC.Range("L38").Value = 5974.9998
$ws_ALC.Range("M38").Value = -326.4000000000001
$ws_ALC.Range("N38").Value = -6718.9998
$ws_ALC.Range("H43").Value = 1256.0769
$ws_ALC.Range("I43").Value = 775.8
$ws_ALC.Range("J43").Value = 1556.25
$ws_ALC.Range("K43").Value = 775.8
$ws_ALC.Range("L43").Value = 1556.25
$ws_ALC.Range("M43").Value = -706.8
$ws_ALC.Range("N43").Value = -1694.25
$ws_ALC.Range("H86").Value = 50024.094
$ws_ALC.Range("I86").Value = 112811.336
$ws_ALC.Range("J86").Value = 2933.6667
$ws_ALC.Range("K86").Value = 112811.336
$ws_ALC.Range("L86").Value = 2933.6667
$ws_ALC.Range("M86").Value = -111688.336
$ws_ALC.Range("N86").Value = -5179.6667
$ws_ALC.Range("H89").Value = 50024.094
$ws_ALC.Range("I89").Value = 112811.336
$ws_ALC.Range("J89").Value = 2933.6667
$ws_ALC.Range("K89").Value = 564056.6799999999
$ws_ALC.Range("L89").Value = 14668.3335
$ws_ALC.Range("M89").Value = -558440.6799999999
$ws_ALC.Range("N89").Value = -25900.3335
$ws_ALC.Range("H100").Value = 1548
$ws_ALC.Range("I100").Value = 1413.3334
$ws_ALC.Range("J100").Value = 1750
$ws_ALC.Range("K100").Value = 1413.3334
$ws_ALC.Range("L100").Value = 1750
$ws_ALC.Range("M100").Value = -872.3334
$ws_ALC.Range("N100").Value = -2832
$ws_ALC.Range("H137").Value = 2990.2
$ws_ALC.Range("I137").Value = 2326.6765
$ws_ALC.Range("K137").Value = 6980.029500000001
$ws_ALC.Range("M137").Value = -4430.029500000001

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 17718
$ws_ARM.Range("I32").Value = 18879.758
$ws_ARM.Range("K32").Value = 18879.758
$ws_ARM.Range("M32").Value = -18592.758
$ws_ARM.Range("H37").Value = 0
$ws_ARM.Range("I37").Value = 0
$ws_ARM.Range("K37").Value = 0
$ws_ARM.Range("M37").ClearContents()
$ws_ARM.Range("H45").Value = 1568.3334
$ws_ARM.Range("I45").Value = 1536.2609
$ws_ARM.Range("J45").Value = 1752.75
$ws_ARM.Range("K45").Value = 1536.2609
$ws_ARM.Range("L45").Value = 1752.75
$ws_ARM.Range("M45").Value = -1159.2609
$ws_ARM.Range("N45").Value = -2506.75
$ws_ARM.Range("H102").Value = 1279729.4
$ws_ARM.Range("I102").Value = 1853987.6
$ws_ARM.Range("K102").Value = 1853987.6
$ws_ARM.Range("M102").Value = -1852365.6
$ws_ARM.Range("H110").Value = 1459.5333
$ws_ARM.Range("I110").Value = 1442.5834
$ws_ARM.Range("J110").Value = 1527.3334
$ws_ARM.Range("K110").Value = 1442.5834
$ws_ARM.Range("L110").Value = 1527.3334
$ws_ARM.Range("M110").Value = 602.4166
$ws_ARM.Range("N110").Value = -5617.3334
$ws_ARM.Range("H122").Value = 2581.3572
$ws_ARM.Range("I122").Value = 2651.2856
$ws_ARM.Range("J122").Value = 2511.4285
$ws_ARM.Range("K122").Value = 7953.8568
$ws_ARM.Range("L122").Value = 7534.2855
$ws_ARM.Range("M122").Value = -5503.8568
$ws_ARM.Range("N122").Value = -12434.2855
$ws_ARM.Range("H132").Value = 5245.881
$ws_ARM.Range("I132").Value = 2430.5652
$ws_ARM.Range("J132").Value = 8653.895
$ws_ARM.Range("K132").Value = 7291.6956
$ws_ARM.Range("L132").Value = 25961.685
$ws_ARM.Range("M132").Value = -4761.6956
$ws_ARM.Range("N132").Value = -31021.685

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 1535.5
$ws_BSM.Range("I20").Value = 1515.4286
$ws_BSM.Range("J20").Value = 1563.6
$ws_BSM.Range("K20").Value = 1515.4286
$ws_BSM.Range("L20").Value = 1563.6
$ws_BSM.Range("M20").Value = -1268.4286
$ws_BSM.Range("N20").Value = -2057.6
$ws_BSM.Range("H107").Value = 1522
$ws_BSM.Range("J107").Value = 1808.6
$ws_BSM.Range("L107").Value = 1808.6
$ws_BSM.Range("N107").Value = -5648.6
$ws_BSM.Range("H132").Value = 56215
$ws_BSM.Range("J132").Value = 65268.75
$ws_BSM.Range("L132").Value = 65268.75
$ws_BSM.Range("N132").Value = -75388.75

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H99").Value = 3673.6
$ws_CRP.Range("I99").Value = 2676.5715
$ws_CRP.Range("J99").Value = 6000
$ws_CRP.Range("K99").Value = 2676.5715
$ws_CRP.Range("L99").Value = 6000
$ws_CRP.Range("M99").Value = -1178.5715
$ws_CRP.Range("N99").Value = -8996
$ws_CRP.Range("H107").Value = 1048.3214
$ws_CRP.Range("I107").Value = 1190.2667
$ws_CRP.Range("J107").Value = 884.53845
$ws_CRP.Range("K107").Value = 1190.2667
$ws_CRP.Range("L107").Value = 884.53845
$ws_CRP.Range("M107").Value = 729.7333000000001
$ws_CRP.Range("N107").Value = -4724.53845
$ws_CRP.Range("H126").Value = 3673.6
$ws_CRP.Range("I126").Value = 2676.5715
$ws_CRP.Range("J126").Value = 6000
$ws_CRP.Range("K126").Value = 8029.7145
$ws_CRP.Range("L126").Value = 18000
$ws_CRP.Range("M126").Value = -5559.7145
$ws_CRP.Range("N126").Value = -22940

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H43").Value = 6000
$ws_CUL.Range("J43").Value = 6000
$ws_CUL.Range("L43").Value = 18000
$ws_CUL.Range("N43").Value = -18228
$ws_CUL.Range("H68").Value = 667.93335
$ws_CUL.Range("J68").Value = 643.25
$ws_CUL.Range("L68").Value = 1929.75
$ws_CUL.Range("N68").Value = -3551.75
$ws_CUL.Range("H71").Value = 667.93335
$ws_CUL.Range("J71").Value = 643.25
$ws_CUL.Range("L71").Value = 5789.25
$ws_CUL.Range("N71").Value = -13901.25
$ws_CUL.Range("H129").Value = 2419.2632
$ws_CUL.Range("I129").Value = 2511.9
$ws_CUL.Range("J129").Value = 2316.3333
$ws_CUL.Range("K129").Value = 7535.700000000001
$ws_CUL.Range("L129").Value = 6948.999899999999
$ws_CUL.Range("M129").Value = -2535.700000000001
$ws_CUL.Range("N129").Value = -16948.9999
$ws_CUL.Range("H132").Value = 1661.1666
$ws_CUL.Range("I132").Value = 1831.5
$ws_CUL.Range("J132").Value = 1490.8334
$ws_CUL.Range("K132").Value = 16483.5
$ws_CUL.Range("L132").Value = 13417.5006
$ws_CUL.Range("M132").Value = -13953.5
$ws_CUL.Range("N132").Value = -18477.5006

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 6076.533
$ws_GSM.Range("I80").Value = 8935.429
$ws_GSM.Range("K80").Value = 8935.429
$ws_GSM.Range("M80").Value = -7937.429
$ws_GSM.Range("H83").Value = 6076.533
$ws_GSM.Range("I83").Value = 8935.429
$ws_GSM.Range("K83").Value = 44677.145
$ws_GSM.Range("M83").Value = -39685.145
$ws_GSM.Range("H107").Value = 428.58823
$ws_GSM.Range("I107").Value = 198.66667
$ws_GSM.Range("K107").Value = 198.66667
$ws_GSM.Range("M107").Value = 1721.33333
$ws_GSM.Range("H132").Value = 4432.375
$ws_GSM.Range("I132").Value = 1734.3823
$ws_GSM.Range("K132").Value = 5203.1469
$ws_GSM.Range("M132").Value = -2673.1469

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 1786.1666
$ws_LTW.Range("I16").Value = 1317
$ws_LTW.Range("K16").Value = 1317
$ws_LTW.Range("M16").Value = -1147
$ws_LTW.Range("H136").Value = 4148.0547
$ws_LTW.Range("I136").Value = 2380.7188
$ws_LTW.Range("K136").Value = 7142.1564
$ws_LTW.Range("M136").Value = -4592.1564

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H21").Value = 0
$ws_WVR.Range("I21").Value = 0
$ws_WVR.Range("K21").Value = 0
$ws_WVR.Range("M21").ClearContents()
$ws_WVR.Range("H35").Value = 0
$ws_WVR.Range("I35").Value = 0
$ws_WVR.Range("K35").Value = 0
$ws_WVR.Range("M35").ClearContents()
